$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.091.88'
$ws.Range("E2").Value = '  +5.48%  '
$ws.Range("D3").Value = '2.382.32'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.67%  '
$ws.Range("D9").Value = '2.379.42'
$ws.Range("E9").Value = '  +3.98%  '
$ws.Range("E10").Value = '  +2.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.58'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.341'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.68%  '
$ws.Range("D15").Value = '2.800.82'
$ws.Range("E15").Value = '  +3.78%  '
$ws.Range("D16").Value = '60.914.65'
$ws.Range("E16").Value = '  +5.22%  '
$ws.Range("E17").Value = '  +2.40%  '
$ws.Range("D18").Value = '2.376.10'
$ws.Range("E18").Value = '  +4.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.46%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '319.90'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.43%  '
$ws.Range("E25").Value = '  +4.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.95%  '
$ws.Range("E28").Value = '  +4.44%  '
$ws.Range("E29").Value = '  +2.86%  '
$ws.Range("D30").Value = '0.0₃0757'
$ws.Range("E30").Value = '  +4.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("E32").Value = '  +7.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.99'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +15.13%  '
$ws.Range("E35").Value = '  +2.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.16'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.75%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.23'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '327.65'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +13.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.59'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '146.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0957'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.35%  '
$ws.Range("E47").Value = '  +1.56%  '
$ws.Range("E48").Value = '  +2.80%  '
$ws.Range("E49").Value = '  +2.43%  '
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("E51").Value = '  +4.69%  '
